# Updated arrival board data (Arribos 141 actualizados - 23)
# Columns: ETA, BANDERA, MIN, ESTADO

$wb = $excel.ActiveWorkbook

$dataTodos = @(
    @("17:12", "16_SANTA ANA", 0, "🚌"),
    @("17:13", "215A_LA PLATA", 1, "🚌"),
    @("17:16", "11_ETCHEVERRY", 4, "🚌"),
    @("17:20", "16_SANTA ANA", 8, "🚌"),
    @("17:20", "26_HERNANDEZ", 8, "📅"),
    @("17:21", "10_OLMOS", 9, "🚌"),
    @("17:28", "14_ABASTO", 16, "🚌"),
    @("17:33", "10_OLMOS", 21, "🚌"),
    @("17:34", "23_HERNANDEZ", 22, "🚌"),
    @("17:36", "27_EL RETIRO", 24, "🚌"),
    @("17:38", "17_ROMERO", 26, "📅"),
    @("17:40", "15_ABASTO", 28, "🚌"),
    @("17:40", "215B_EL PATO", 28, "📅"),
    @("17:45", "15_ABASTO", 33, "🚌"),
    @("17:50", "16_P MOR-167 Y 521", 38, "🚌"),
    @("17:52", "81_EL PELIGRO", 40, "📅"),
    @("18:03", "215C_LA PLATA", 51, "🚌"),
    @("18:04", "17_ROMERO", 52, "🚌"),
    @("18:08", "14_ABASTO", 56, "🚌"),
    @("18:12", "23_HERNANDEZ", 60, "🚌"),
    @("18:20", "26_HERNANDEZ", 68, "🚌"),
    @("18:27", "215C_EL PATO", 75, "🚌"),
    @("18:31", "11X44_ETCHEVERRY", 79, "🚌"),
    @("18:39", "84_COLONIA URQUIZA-ESC 49", 87, "🚌"),
    @("18:47", "14X44_ABASTO", 95, "🚌"),
    @("18:51", "215A_LA PLATA", 99, "🚌"),
    @("18:58", "215A_EL PATO", 106, "📅"),
    @("19:03", "215B_LP-P MOR-1 Y 57", 111, "📅"),
    @("19:04", "11_ETCHEVERRY", 112, "🚌"),
    @("19:10", "16_P MOR-SANTA ANA", 118, "📅")
)

$data215 = @(
    @("17:13", "215A_LA PLATA", 1, "🚌"),
    @("17:40", "215B_EL PATO", 28, "📅"),
    @("18:03", "215C_LA PLATA", 51, "🚌"),
    @("18:27", "215C_EL PATO", 75, "🚌"),
    @("18:51", "215A_LA PLATA", 99, "🚌"),
    @("18:58", "215A_EL PATO", 106, "📅"),
    @("19:03", "215B_LP-P MOR-1 Y 57", 111, "📅")
)

$wsTodos = $wb.Worksheets.Item("TODOS")
for ($i = 0; $i -lt $dataTodos.Length; $i++) {
    $row = $dataTodos[$i]
    $r = $i + 2
    $wsTodos.Cells.Item($r, 1).Value = $row[0]
    $wsTodos.Cells.Item($r, 2).Value = $row[1]
    $wsTodos.Cells.Item($r, 3).Value = $row[2]
    $wsTodos.Cells.Item($r, 4).Value = $row[3]
}

$ws215 = $wb.Worksheets.Item("215")
for ($i = 0; $i -lt $data215.Length; $i++) {
    $row = $data215[$i]
    $r = $i + 2
    $ws215.Cells.Item($r, 1).Value = $row[0]
    $ws215.Cells.Item($r, 2).Value = $row[1]
    $ws215.Cells.Item($r, 3).Value = $row[2]
    $ws215.Cells.Item($r, 4).Value = $row[3]
}

$wsCombinadas = $wb.Worksheets.Item("COMBINADAS")
for ($i = 0; $i -lt $dataTodos.Length; $i++) {
    $row = $dataTodos[$i]
    $r = $i + 2
    $wsCombinadas.Cells.Item($r, 1).Value = $row[0]
    $wsCombinadas.Cells.Item($r, 2).Value = $row[1]
    $wsCombinadas.Cells.Item($r, 3).Value = $row[2]
    $wsCombinadas.Cells.Item($r, 4).Value = $row[3]
}
